# no-op test
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$t = $s.Shapes.Item(1).TextFrame.TextRange.Text
